# Add a new "TaxonRelation" worksheet right after "NameRelation", carrying
# the same look & feel as the other relation-style sheets (CoLDP template),
# and populate its header row with the new taxon-concept/species-interaction
# relation columns.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("NameRelation")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "TaxonRelation"

# Header row: taxonID | relatedTaxonID | type | referenceID | remarks
$newSheet.Range("A1").Value = "taxonID"
$newSheet.Range("B1").Value = "relatedTaxonID"
$newSheet.Range("C1").Value = "type"
$newSheet.Range("D1").Value = "referenceID"
$newSheet.Range("E1").Value = "remarks"

# Match the wide column layout used by the sibling relation sheets.
$newSheet.Columns("A:E").ColumnWidth = 45.33

# Match the zoomed-in view used across the workbook's relation sheets.
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 160
